$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.263.97"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.592.33"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.17"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.246"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.98"
$ws.Range("E10").Value = "  -2.00%  "

$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.815.90"
$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.594.38"
$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("E14").Value = "  -1.08%  "

$ws.Range("E15").Value = "  -2.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.82"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.249.09"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.06"
$ws.Range("E19").Value = "  +1.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -1.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("E23").Value = "  +0.63%  "

$ws.Range("E24").Value = "  -3.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.68"

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.97"
$ws.Range("E27").Value = "  -1.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.113"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("E30").Value = "  -2.09%  "

$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.409.64"
$ws.Range("E33").Value = "  +5.35%  "

$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("E36").Value = "  -1.42%  "

$ws.Range("E37").Value = "  -3.96%  "

$ws.Range("E38").Value = "  -1.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.822"
$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.77"
$ws.Range("E40").Value = "  -0.62%  "

$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.954"
$ws.Range("E42").Value = "  -8.11%  "

$ws.Range("E43").Value = "  +1.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.761"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.728.06"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.92"
$ws.Range("E46").Value = "  -1.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.06"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0501"
$ws.Range("E49").Value = "  -0.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0951"
$ws.Range("E50").Value = "  -3.24%  "

$ws.Range("E51").Value = "  +0.04%  "
